$d = $word.ActiveDocument
$wordMl = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$p = $d.Paragraphs.Item(1)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:r><w:t>ContosoLearn Market Research</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(2)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AdatumLearn: AdatumLearn is a top AI-powered learning platform that uses artificial intelligence to enrich eLearning with features that automate a variety of tasks. It is known for its content authoring capabilities and adaptive learning technology.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(3)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AdventureLearn: AdventureLearn is another AI-powered learning platform that offers personalized learning experiences and data-driven recommendations.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(4)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>AlpineTraining: AlpineTraining is a mobile-first learning platform that focuses on microlearning.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(5)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>Bellows OnDemand: Bellows OnDemand is a comprehensive learning solution that offers content creation and social collaboration.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(6)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>FabrikamLearning: FabrikamLearning provides a suite of learning platforms that cater to different learning needs.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(7)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>FirstUp Cards: FirstUp Cards is a mobile learning app that is ideal for training on safety procedures, compliance, new product knowledge or any other type of training scenario.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(8)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>Munson&#x27;sLearn: Munson&#x27;sLearn is designed to enable businesses to train their employees, partners, and customers.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(9)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>LibertyLearn: LibertyLearn is a fast LMS for your mission-critical project.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(10)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">WoodgroveLMS: WoodgroveLMS is a functional and attractive learning management system built to provide </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a best</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>-in-class training experience.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(11)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>NorthwindWorlds: NorthwindWorlds is a powerful, easy-to-use, and reliable training solution for individuals and enterprises.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(12)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>ProsewareLearn: ProsewareLearn is an online education company that offers a variety of video training courses for software developers, IT administrators, and creative professionals through its website.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(13)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>RelecloudLearn: RelecloudLearn is an American online learning platform that offers massive open online courses (MOOC), specializations, and degrees in a variety of subjects.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(14)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="14"/></w:numPr></w:pPr><w:r><w:t>TreyAcademy: TreyAcademy is an online learning platform aimed at professional adults and students, developed in May 2010.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

$p = $d.Paragraphs.Item(15)
$r = $p.Range
$xml = '<?xml version=''1.0''?><pkg:package xmlns:pkg=''http://schemas.microsoft.com/office/2006/xmlPackage''><pkg:part pkg:name=''/word/document.xml'' pkg:contentType=''application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml''><pkg:xmlData><w:document xmlns:w=''http://schemas.openxmlformats.org/wordprocessingml/2006/main''><w:body><w:p><w:r><w:t xml:space="preserve">These platforms have a significant market presence and are widely recognized for their AI-powered features, such as personalized learning experiences, data-driven recommendations, and automation of tasks. They are transforming the eLearning landscape by leveraging AI to deliver more engaging, rewarding, and personalized learning experiences. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
[void]$r.InsertXML($xml)

